# example/exampledata.xlsx — "further finished, some problem with times in
# the example sheet"
#
# Adds a "format" column to the events sheet (marking the keynote and the
# workshop slots as PLENARY events) and an "items" column to the sessions
# sheet, then leaves the cursor on the sessions sheet (cell D2), which is
# where the author ended up after finishing the edit.

$wb = $excel.ActiveWorkbook

$wsInstructions = $wb.Worksheets.Item("instructions")
$wsItems        = $wb.Worksheets.Item("items")
$wsSessions     = $wb.Worksheets.Item("sessions")
$wsEvents       = $wb.Worksheets.Item("events")
$wsAuthors      = $wb.Worksheets.Item("authors")

# --- content edits -------------------------------------------------------

# sessions: new trailing column "items"
$wsSessions.Range("G1").Value = "items"

# events: new trailing column "format", marking the two plenary slots
$wsEvents.Range("G1").Value = "format"
$wsEvents.Range("G2").Value = "PLENARY"
$wsEvents.Range("G6").Value = "PLENARY"

# --- restore each sheet's own resting selection ---------------------------

$wsInstructions.Select()
$wsInstructions.Range("A6").Select()

$wsItems.Select()
$wsItems.Range("D31").Select()

$wsEvents.Select()
$wsEvents.Range("F4").Select()

$wsAuthors.Select()
$wsAuthors.Range("C3").Select()

# sessions ends up the active sheet, cursor on D2
$wsSessions.Select()
$wsSessions.Range("D2").Select()
